# Apply cryptocurrency price/volume/hour refresh for "Updated symbol list" commit.
# Numeric-looking values are written with a leading apostrophe so Excel keeps them
# as text (matching the original t="inlineStr" cell type) instead of auto-converting
# them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'276.13"
$ws.Range("G2").Value = "'14"

# Row 3
$ws.Range("D3").Value = "'21.02"
$ws.Range("G3").Value = "'14"

# Row 4
$ws.Range("D4").Value = "'6.212"
$ws.Range("G4").Value = "'14"

# Row 5
$ws.Range("D5").Value = "'0.06182"
$ws.Range("G5").Value = "'14"

# Row 6
$ws.Range("D6").Value = "'3.582"
$ws.Range("G6").Value = "'14"

# Row 7
$ws.Range("D7").Value = "'6.568"
$ws.Range("G7").Value = "'14"

# Row 8
$ws.Range("D8").Value = "'1.493"
$ws.Range("G8").Value = "'14"

# Row 9
$ws.Range("D9").Value = "'0.8232"
$ws.Range("G9").Value = "'14"

# Row 10
$ws.Range("D10").Value = "'0.01379"
$ws.Range("G10").Value = "'14"

# Row 11
$ws.Range("D11").Value = "'0.1616"
$ws.Range("G11").Value = "'14"

# Row 12
$ws.Range("D12").Value = "'0.08239"
$ws.Range("G12").Value = "'14"

# Row 13
$ws.Range("D13").Value = "'0.03512"
$ws.Range("G13").Value = "'14"

# Row 14
$ws.Range("D14").Value = "'0.03101"
$ws.Range("G14").Value = "'14"

# Row 15
$ws.Range("D15").Value = "'0.09124"
$ws.Range("G15").Value = "'14"

# Row 16
$ws.Range("D16").Value = "'3.764"
$ws.Range("G16").Value = "'14"

# Row 17
$ws.Range("D17").Value = "'0.001625"
$ws.Range("G17").Value = "'14"

# Row 18
$ws.Range("D18").Value = "'0.04690"
$ws.Range("G18").Value = "'14"

# Row 19
$ws.Range("D19").Value = "'0.006437"
$ws.Range("G19").Value = "'14"

# Row 20
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "'0.006162"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("G20").Value = "'14"

# Row 21
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "'0.001068"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("G21").Value = "'14"

# Row 22
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.0001502"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("G22").Value = "'14"

# Row 23
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.803"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("G23").Value = "'14"

# Row 24
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.282"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("G24").Value = "'14"

# Row 25
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D25").Value = "'0.3389"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"
$ws.Range("G25").Value = "'14"

# Row 26
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D26").Value = "'0.1222"
$ws.Range("E26").Value = "25ProBitTokenPROB"
$ws.Range("G26").Value = "'14"

# Row 27
$ws.Range("B27").Value = "AAXToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab"
$ws.Range("D27").Value = "'0.3999"
$ws.Range("E27").Value = "26AAXTokenAAB"
$ws.Range("G27").Value = "'14"

# Row 28
$ws.Range("D28").Value = "'0.0001583"
$ws.Range("G28").Value = "'14"

# Row 29
$ws.Range("G29").Value = "'14"

# Row 30
$ws.Range("G30").Value = "'14"

# Row 31
$ws.Range("G31").Value = "'14"

# Row 32
$ws.Range("G32").Value = "'14"

# Row 33
$ws.Range("G33").Value = "'14"

# Row 34
$ws.Range("G34").Value = "'14"

# Row 35
$ws.Range("G35").Value = "'14"

# Row 36
$ws.Range("G36").Value = "'14"

# Row 37
$ws.Range("G37").Value = "'14"

# Row 38
$ws.Range("G38").Value = "'14"

# Row 39
$ws.Range("G39").Value = "'14"

# Row 40
$ws.Range("D40").Value = "'0.04652"
$ws.Range("G40").Value = "'14"

# Row 41
$ws.Range("D41").Value = "'0.007065"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("G41").Value = "'14"

# Row 42
$ws.Range("D42").Value = "'0.004606"
$ws.Range("E42").Value = "41CEJICEJIBestin24h"
$ws.Range("G42").Value = "'14"

# Row 43
$ws.Range("D43").Value = "'0.1101"
$ws.Range("G43").Value = "'14"

# Row 44
$ws.Range("D44").Value = "'0.01083"
$ws.Range("G44").Value = "'14"

# Row 45
$ws.Range("D45").Value = "'0.00006168"
$ws.Range("G45").Value = "'14"

# Row 46
$ws.Range("G46").Value = "'14"

# Row 47
$ws.Range("D47").Value = "'0.8467"
$ws.Range("G47").Value = "'14"

# Row 48
$ws.Range("D48").Value = "'0.002619"
$ws.Range("G48").Value = "'14"

# Row 49
$ws.Range("D49").Value = "'0.00001903"
$ws.Range("G49").Value = "'14"

# Row 50
$ws.Range("D50").Value = "'0.01242"
$ws.Range("G50").Value = "'14"

# Row 51
$ws.Range("G51").Value = "'14"
